# Add lab 2 demographic data for subjects 4-6 (rows 5-7): gender + age.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "f"
$ws.Range("C5").Value = 21

$ws.Range("B6").Value = "f"
$ws.Range("C6").Value = 19

$ws.Range("B7").Value = "m"
$ws.Range("C7").Value = 24

# Leave the selection where the author's last edit landed.
$ws.Range("C8").Select()
